$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Locate the "KEY ACHIEVEMENTS AND IMPACT" section bounds so the edits below
# are scoped to that section only (near-duplicate bullet text also appears
# earlier, under the Siege Analytics experience entry, and must be left
# untouched).
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text.Trim()
    if ($t -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $startIdx = $i
    }
    if ($startIdx -ge 1 -and $i -gt $startIdx -and $t -eq "TECHNICAL SKILLS") {
        $endIdx = $i
        break
    }
}

$sectionStart = $paras.Item($startIdx).Range.Start
$sectionEnd = $paras.Item($endIdx).Range.Start
$sectionRange = $d.Range($sectionStart, $sectionEnd)

# 1) Rewrite the first four bullet paragraphs as impact-focused statements.
$sectionRange.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions",
    2) | Out-Null

$sectionRange.Find.Execute(
    "Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "178% accuracy improvement in racial classification algorithms",
    2) | Out-Null

$sectionRange.Find.Execute(
    "Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%",
    2) | Out-Null

$sectionRange.Find.Execute(
    "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "`$4.7M savings enabled nonprofit access",
    2) | Out-Null

# 2) Remove the "Built real-time FEC analysis systems..." bullet paragraph
#    entirely (it has no replacement in the rewritten list). Iterate the
#    document's paragraph collection directly (bounded by the section
#    indices found above) rather than Range.Paragraphs.
$target = $null
for ($i = $startIdx; $i -lt $endIdx; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Built real-time FEC analysis systems*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# 3) Replace the final remaining bullet paragraph's text.
$sectionRange.Find.Execute(
    "Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations",
    2) | Out-Null
